# Refresh the "cryptos" price/volume snapshot (scheduled GitHub Actions update).
# Price cells that look like plain decimals (e.g. "0.999") are written with a
# leading apostrophe so Excel stores them as text, matching the sheet's
# existing inline-string convention instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.424.52"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "2.512.49"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'588.77"
$ws.Range("E5").Value = "  -0.55%  "

$ws.Range("D6").Value = "'169.65"
$ws.Range("E6").Value = "  -3.53%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  -2.32%  "

$ws.Range("D9").Value = "2.509.12"
$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("E10").Value = "  -3.48%  "

$ws.Range("D11").Value = "'0.166"
$ws.Range("E11").Value = "  +0.98%  "

$ws.Range("E12").Value = "  -1.43%  "

$ws.Range("E13").Value = "  -3.80%  "

$ws.Range("D14").Value = "3.001.70"
$ws.Range("E14").Value = "  +0.78%  "

$ws.Range("D15").Value = "'26.02"
$ws.Range("E15").Value = "  -2.93%  "

$ws.Range("E16").Value = "  -2.78%  "

$ws.Range("D17").Value = "67.302.22"
$ws.Range("E17").Value = "  -0.46%  "

$ws.Range("D18").Value = "2.540.92"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'11.71"
$ws.Range("E19").Value = "  +2.34%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'8.02"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("D21").Value = "'364.29"
$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("B22").Value = "Binance-PegBSC-USD"
$ws.Range("C22").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D22").Value = "'1.55"
$ws.Range("E22").Value = "  +55.28%  "

$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D23").Value = "'4.08"
$ws.Range("E23").Value = "  -2.59%  "

$ws.Range("D24").Value = "'4.49"
$ws.Range("E24").Value = "  -3.36%  "

$ws.Range("D25").Value = "'71.70"
$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("E27").Value = "  -6.64%  "

$ws.Range("D28").Value = "'9.60"
$ws.Range("E28").Value = "  -6.58%  "

$ws.Range("D29").Value = "2.650.59"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("D30").Value = "0.0₃0937"
$ws.Range("E30").Value = "  -5.45%  "

$ws.Range("D31").Value = "'528.38"
$ws.Range("E31").Value = "  -2.54%  "

$ws.Range("E32").Value = "  -0.73%  "

$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("E34").Value = "  -4.99%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("E36").Value = "  -1.82%  "

$ws.Range("D37").Value = "'157.96"
$ws.Range("E37").Value = "  +1.37%  "

$ws.Range("D38").Value = "'19.16"
$ws.Range("E38").Value = "  +1.83%  "

$ws.Range("E39").Value = "  -2.89%  "

$ws.Range("D40").Value = "'18.59"
$ws.Range("E40").Value = "  -0.37%  "

$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("E42").Value = "  -4.27%  "

$ws.Range("E43").Value = "  -3.34%  "

$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("E45").Value = "  -3.99%  "

$ws.Range("D46").Value = "'39.34"

$ws.Range("D47").Value = "'146.82"
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").Value = "'0.542"
$ws.Range("E49").Value = "  -3.36%  "

$ws.Range("D50").Value = "0.0₆0269"
$ws.Range("E50").Value = "  -3.98%  "

$ws.Range("E51").Value = "  -0.18%  "

